# Regenerate column G ("K" = strikeouts) using actual strikeout totals (K)
# instead of the previous proxy "Strike#" value, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$kValues = @{
    2 = 0
    3 = 1
    4 = 1
    5 = 1
    6 = 0
    7 = 1
    8 = 1
    9 = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 2
    30 = 3
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 3
    38 = 1
    39 = 0
    40 = 1
    41 = 1
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 3
    47 = 1
    48 = 2
    49 = 0
    50 = 0
    51 = 0
    52 = 2
    53 = 3
    54 = 0
    55 = 3
    56 = 1
    57 = 1
    58 = 2
    59 = 0
    60 = 0
    61 = 0
    62 = 2
    63 = 2
    64 = 0
    65 = 1
    66 = 1
    67 = 2
    68 = 0
    69 = 1
    70 = 3
    71 = 1
    72 = 0
    73 = 1
    74 = 2
    75 = 2
    76 = 1
    77 = 2
    78 = 0
    79 = 2
    80 = 0
    81 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Keys.Count) K values in column G"
